$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1026.3889
$ws.Range("I15").Value = 1026.3889
$ws.Range("K15").Value = 3079.1667
$ws.Range("M15").Value = -2910.1667
$ws.Range("H93").Value = 83996
$ws.Range("J93").Value = 83996
$ws.Range("L93").Value = 83996
$ws.Range("N93").Value = -88988
$ws.Range("H112").Value = 6987.375
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("H125").Value = 1375
$ws.Range("I125").Value = 750
$ws.Range("J125").Value = 3250
$ws.Range("K125").Value = 6750
$ws.Range("L125").Value = 29250
$ws.Range("M125").Value = -4290
$ws.Range("N125").Value = -34170
$ws.Range("H129").Value = 913.0909
$ws.Range("J129").Value = 880.7843
$ws.Range("L129").Value = 2642.3529
$ws.Range("N129").Value = -12642.3529
$ws.Range("H132").Value = 1259.5
$ws.Range("I132").Value = 1273.1578
$ws.Range("K132").Value = 3819.4734
$ws.Range("M132").Value = -1289.4734
$ws.Range("H137").Value = 1564.7059
$ws.Range("I137").Value = 1373.3334
$ws.Range("K137").Value = 4120.0002
$ws.Range("M137").Value = -1570.0002
$ws.Range("M112").ClearContents()

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 50668.668
$ws.Range("J3").Value = 50668.668
$ws.Range("L3").Value = 50668.668
$ws.Range("N3").Value = -50898.668
$ws.Range("H32").Value = 4021.2642
$ws.Range("I32").Value = 2286.2092
$ws.Range("K32").Value = 2286.2092
$ws.Range("M32").Value = -1999.2092
$ws.Range("H45").Value = 1777
$ws.Range("I45").Value = 1350
$ws.Range("K45").Value = 1350
$ws.Range("M45").Value = -973

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 42580.6
$ws.Range("I11").Value = 1444
$ws.Range("K11").Value = 1444
$ws.Range("M11").Value = -1304
$ws.Range("H26").Value = 39997.5
$ws.Range("I26").Value = 39997.5
$ws.Range("K26").Value = 39997.5
$ws.Range("M26").Value = -39705.5
$ws.Range("H96").Value = 40000
$ws.Range("I96").Value = 40000
$ws.Range("K96").Value = 40000
$ws.Range("M96").Value = -37254
$ws.Range("H122").Value = 68000
$ws.Range("J122").Value = 68000
$ws.Range("L122").Value = 68000
$ws.Range("N122").Value = -77800

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2594.3865
$ws.Range("I31").Value = 1654.1034
$ws.Range("J31").Value = 4412.2666
$ws.Range("K31").Value = 1654.1034
$ws.Range("L31").Value = 4412.2666
$ws.Range("M31").Value = -1359.1034
$ws.Range("N31").Value = -5002.2666
$ws.Range("H34").Value = 2594.3865
$ws.Range("I34").Value = 1654.1034
$ws.Range("J34").Value = 4412.2666
$ws.Range("K34").Value = 1654.1034
$ws.Range("L34").Value = 4412.2666
$ws.Range("M34").Value = -1452.1034
$ws.Range("N34").Value = -4816.2666
$ws.Range("H105").Value = 1161.7778
$ws.Range("I105").Value = 1055.625
$ws.Range("K105").Value = 1055.625
$ws.Range("M105").Value = 691.375
$ws.Range("H132").Value = 1931.3334
$ws.Range("I132").Value = 1078.84
$ws.Range("J132").Value = 3868.818
$ws.Range("K132").Value = 3236.52
$ws.Range("L132").Value = 11606.454
$ws.Range("M132").Value = -706.5199999999995
$ws.Range("N132").Value = -16666.454

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 183.57143
$ws.Range("I2").Value = 246.25
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 1477.5
$ws.Range("L2").Value = 600
$ws.Range("M2").Value = -1364.5
$ws.Range("N2").Value = -826
$ws.Range("H5").Value = 849.125
$ws.Range("I5").Value = 734.6667
$ws.Range("J5").Value = 917.8
$ws.Range("K5").Value = 2204.0001
$ws.Range("L5").Value = 2753.4
$ws.Range("M5").Value = -2092.0001
$ws.Range("N5").Value = -2977.4
$ws.Range("H11").Value = 846
$ws.Range("I11").Value = 525
$ws.Range("J11").Value = 1488
$ws.Range("K11").Value = 1575
$ws.Range("L11").Value = 4464
$ws.Range("M11").Value = -1435
$ws.Range("N11").Value = -4744
$ws.Range("I26").Value = 300
$ws.Range("J26").Value = 326
$ws.Range("K26").Value = 900
$ws.Range("L26").Value = 978
$ws.Range("M26").Value = -612
$ws.Range("N26").Value = -1554
$ws.Range("H113").Value = 7436.3335
$ws.Range("J113").Value = 874
$ws.Range("L113").Value = 2622
$ws.Range("N113").Value = -6962
$ws.Range("H122").Value = 879.9167
$ws.Range("I122").Value = 694.75
$ws.Range("J122").Value = 972.5
$ws.Range("K122").Value = 6252.75
$ws.Range("L122").Value = 8752.5
$ws.Range("M122").Value = -3802.75
$ws.Range("N122").Value = -13652.5
$ws.Range("H131").Value = 15588.596
$ws.Range("J131").Value = 17776.56
$ws.Range("L131").Value = 53329.68000000001
$ws.Range("N131").Value = -63409.68000000001
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("H135").Value = 849.125
$ws.Range("I135").Value = 734.6667
$ws.Range("J135").Value = 917.8
$ws.Range("K135").Value = 6612.0003
$ws.Range("L135").Value = 8260.199999999999
$ws.Range("M135").Value = -4077.0003
$ws.Range("N135").Value = -13330.2
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2266266.8
$ws.Range("J132").Value = 11563
$ws.Range("L132").Value = 34689
$ws.Range("N132").Value = -39749

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2499.7
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2499.7
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2499.7
$ws.Range("N46").Value = -2875.7
$ws.Range("H82").Value = 3378.4
$ws.Range("I82").Value = 1300.6666
$ws.Range("K82").Value = 1300.6666
$ws.Range("M82").Value = -939.6666
$ws.Range("H85").Value = 3378.4
$ws.Range("I85").Value = 1300.6666
$ws.Range("K85").Value = 1300.6666
$ws.Range("M85").Value = -52.66660000000002
$ws.Range("H93").Value = 424.3
$ws.Range("I93").Value = 305.375
$ws.Range("K93").Value = 305.375
$ws.Range("M93").Value = 942.625
$ws.Range("M46").ClearContents()

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2419.4
$ws.Range("I62").Value = 2232.6667
$ws.Range("J62").Value = 2699.5
$ws.Range("K62").Value = 2232.6667
$ws.Range("L62").Value = 2699.5
$ws.Range("M62").Value = -1608.6667
$ws.Range("N62").Value = -3947.5
$ws.Range("H65").Value = 2419.4
$ws.Range("I65").Value = 2232.6667
$ws.Range("J65").Value = 2699.5
$ws.Range("K65").Value = 11163.3335
$ws.Range("L65").Value = 13497.5
$ws.Range("M65").Value = -8043.333500000001
$ws.Range("N65").Value = -19737.5
$ws.Range("H122").Value = 118828.125
$ws.Range("I122").Value = 145695
$ws.Range("J122").Value = 2405
$ws.Range("K122").Value = 437085
$ws.Range("L122").Value = 7215
$ws.Range("M122").Value = -434635
$ws.Range("N122").Value = -12115
$ws.Range("H132").Value = 2158.6667
$ws.Range("I132").Value = 1653.8235
$ws.Range("J132").Value = 3384.7144
$ws.Range("K132").Value = 4961.470499999999
$ws.Range("L132").Value = 10154.1432
$ws.Range("M132").Value = -2431.470499999999
$ws.Range("N132").Value = -15214.1432
$ws.Range("H139").Value = 59955
$ws.Range("J139").Value = 59955
$ws.Range("L139").Value = 59955
$ws.Range("N139").Value = -70235
